$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.749.89'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.042.44'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'227.47"
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('D7').Value = "'59.95"
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('E10').Value = '  +3.17%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '2.345.25'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = "'14.33"
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').Value = "'21.22"
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = "'5.45"
$ws.Range('E15').Value = '  +4.93%  '
$ws.Range('D16').Value = "'0.762"
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '2.027.32'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '37.723.60'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = "'5.93"
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').Value = "'69.35"
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').Value = "'223.48"
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').Value = "'168.77"
$ws.Range('E26').Value = '  +2.00%  '
$ws.Range('D27').Value = "'9.33"
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').Value = "'18.75"
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'4.48"
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = "'6.52"
$ws.Range('E36').Value = '  +2.62%  '
$ws.Range('E37').Value = '  +4.04%  '
$ws.Range('E38').Value = '  +6.03%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = "'17.93"
$ws.Range('E40').Value = '  +6.27%  '
$ws.Range('D41').Value = '1.535.76'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = "'97.83"
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('D45').Value = "'0.0902"
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('D46').Value = "'4.14"
$ws.Range('E46').Value = '  +5.71%  '
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = "'7.01"
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').Value = '2.233.89'
